$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns D..R ---
$ws.Range("D1").Value = "Camera On While Delivering"
$ws.Range("E1").Value = "Class Started on Time"
$ws.Range("F1").Value = "Zoom Poll Taken / Feedback Poll Taken"
$ws.Range("G1").Value = "Total attendees (online + offline)"
$ws.Range("H1").Value = "Resolution of Non Tech query"
$ws.Range("I1").Value = "Resolution of Tech query"
$ws.Range("J1").Value = "Refer and earn slide shown"
$ws.Range("K1").Value = "Participant Engagement"
$ws.Range("L1").Value = "Technical glitch (if any)"
$ws.Range("M1").Value = "Was there any disruption during the session?"
$ws.Range("N1").Value = "Comments"
$ws.Range("O1").Value = "2025-04-06 18:45:35"
$ws.Range("P1").Value = "2025-04-06 18:46:29"
$ws.Range("Q1").Value = "2025-04-06 18:47:54"
$ws.Range("R1").Value = "2025-04-06 18:49:22"

# Apply the same header formatting (bold, centered, bordered) used by A1:C1 to
# the newly-added header cells, by copy/pasting the format from A1 (this keeps
# the style index stable instead of minting a subtly-different duplicate xf).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:R1").PasteSpecial(-4122) | Out-Null

# --- Row 2 ---
$ws.Range("B2").Value = "15-FEB-25-CDS-BUN-021-WEM0930-BAN & 28-Dec-24-CDS-BUN-021-WEM09-BAN"
$ws.Range("C2").Value = "29-MAR-25-MLE-113-WEM09-BAN (CONTINUE)"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "Yes"
$ws.Range("I2").Value = "Yes"
$ws.Range("J2").Value = "Yes"
$ws.Range("K2").Value = "Yes"
$ws.Range("L2").Value = "No"
$ws.Range("M2").Value = "No"
$ws.Range("N2").Value = "test"
$ws.Range("O2").Value = "Error: name 'tk' is not defined"
$ws.Range("P2").Value = "Submitted"
$ws.Range("Q2").Value = "Submitted"
$ws.Range("R2").Value = "Submitted"

# --- Row 3 ---
$ws.Range("B3").Value = "test"
$ws.Range("C3").Value = "test"
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = 43
$ws.Range("H3").Value = "Yes"
$ws.Range("I3").Value = "Yes"
$ws.Range("J3").Value = "Yes"
$ws.Range("K3").Value = "Yes"
$ws.Range("L3").Value = "No"
$ws.Range("M3").Value = "No"
$ws.Range("N3").Value = "test 2"
$ws.Range("O3").Value = "Error: name 'tk' is not defined"
$ws.Range("P3").Value = "Submitted"
$ws.Range("Q3").Value = "Submitted"
$ws.Range("R3").Value = "Submitted"
